$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.159667134284973
$ws.Range("B1").Value = 3.264880418777466
$ws.Range("C1").Value = 5.288897037506104
$ws.Range("D1").Value = 2.25877571105957
$ws.Range("E1").Value = 1.376042723655701
